# Replace the multiplication problems in the document's table cells.
# Each old expression is unique within the document, so a simple
# Find/Replace (MatchWholeWord, not wildcard) for each pair is safe.

$d = $word.ActiveDocument

$pairs = @(
    @("12×60=", "45×55="),
    @("60×20=", "41×14="),
    @("22×81=", "40×81="),
    @("54×87=", "89×29="),
    @("60×94=", "50×11="),
    @("61×15=", "63×46="),
    @("68×16=", "37×52="),
    @("29×11=", "99×29="),
    @("85×88=", "99×75="),
    @("82×89=", "16×65="),
    @("88×63=", "52×99="),
    @("49×77=", "11×97="),
    @("21×66=", "97×27="),
    @("15×67=", "37×93="),
    @("62×36=", "64×50="),
    @("48×99=", "28×59="),
    @("12×75=", "45×46="),
    @("19×79=", "35×74="),
    @("86×32=", "21×52="),
    @("61×73=", "35×96="),
    @("89×85=", "57×38="),
    @("82×53=", "52×36="),
    @("73×85=", "72×60="),
    @("27×22=", "65×18="),
    @("29×61=", "76×55=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
